$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42605.885787037034
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("B3").Value = -6
$ws.Range("C3").Value = 61
$ws.Range("D3").Value = 37
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = 71
$ws.Range("G3").Value = 12658
$ws.Range("H3").Value = 6725
$ws.Range("I3").Value = 876
$ws.Range("J3").Value = 75
$ws.Range("K3").Value = 46
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 10
$ws.Range("N3").Value = "Bag"
